$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new "Tier" column (C2:C7) added by this commit.
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 9
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 3
$ws.Range("C7").Value = 7

# Move the active selection, matching the author's final cursor position.
$ws.Range("D10").Select()
